$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "slowo"
$ws.Range("A3").Value = "ziemniak"
$ws.Range("B2").Value = "word"
$ws.Range("B3").Value = "potato"
